$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# 1. Rewrite the table data (now A:E instead of A:D, 8 rows)
#    Column B is written with a placeholder everywhere first so
#    formatting can be applied uniformly, then the genuinely blank
#    rows (4 and 7) are cleared completely afterwards.
# ---------------------------------------------------------------
$rows = @(
    @("Ref No",       "QC Require", "Status",      "Supplier Cost Value($)", "Sunrise Value($)"),
    @("HIYU00040082", "Yes",        "Hold",        100,                      200),
    @("HIYU00040084", "Yes",        "Confirm",     200,                      300),
    @("HIYU00040083", $null,        "Confirm",     300,                      400),
    @("HIYU00040086", "No",         "QC Pending",  400,                      500),
    @("HIYU00040087", "No",         "QC Reject",   500,                      600),
    @("HIYO00040088", $null,        "Confirm",     600,                      700),
    @("HIYU00040088", "Yes",        "Bidded",      700,                      800)
)

for ($r = 0; $r -lt $rows.Length; $r++) {
    $rowNum = $r + 1
    $vals = $rows[$r]

    $ws.Cells.Item($rowNum, 1).Value = $vals[0]
    $ws.Cells.Item($rowNum, 2).Value = "placeholder"
    $ws.Cells.Item($rowNum, 3).Value = $vals[2]
    $ws.Cells.Item($rowNum, 4).Value = $vals[3]
    $ws.Cells.Item($rowNum, 5).Value = $vals[4]
}

# ---------------------------------------------------------------
# 2. Formatting: columns A-C centered text (reuse the format that
#    column A/B already carry), columns D-E centered numbers with
#    #,##0.00 (reuse the format column D already carries).
# ---------------------------------------------------------------
$ws.Range("A1").Copy()
$ws.Range("C1:C8").PasteSpecial(-4122)
$ws.Range("A1:B8").Copy()
$ws.Range("A1:B8").PasteSpecial(-4122)

$ws.Range("D1:D8").Copy()
$ws.Range("D1:E8").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ---------------------------------------------------------------
# 3. Now that styles are settled, write the real column B values.
#    Rows 4 and 7 are cleared completely so no <c> element remains.
# ---------------------------------------------------------------
for ($r = 0; $r -lt $rows.Length; $r++) {
    $rowNum = $r + 1
    $vals = $rows[$r]
    if ($vals[1] -eq $null) {
        $ws.Cells.Item($rowNum, 2).Clear()
    } else {
        $ws.Cells.Item($rowNum, 2).Value = $vals[1]
    }
}

# Re-assert the data values/number format on D:E (PasteSpecial above
# only copied formatting, but make sure values are the right numbers).
for ($r = 0; $r -lt $rows.Length; $r++) {
    $rowNum = $r + 1
    $vals = $rows[$r]
    $ws.Cells.Item($rowNum, 4).Value = $vals[3]
    $ws.Cells.Item($rowNum, 5).Value = $vals[4]
}

# ---------------------------------------------------------------
# 4. Column widths
# ---------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 15.7109375
$ws.Columns.Item(2).ColumnWidth = 15.5703125
$ws.Columns.Item(3).ColumnWidth = 11.28515625
$ws.Columns.Item(4).ColumnWidth = 25.5703125
$ws.Columns.Item(5).ColumnWidth = 20.28515625

# ---------------------------------------------------------------
# 5. AutoFilter over the new A1:E1 header range
# ---------------------------------------------------------------
if ($ws.AutoFilterMode) {
    $ws.AutoFilterMode = $false
}
$ws.Range("A1:E1").AutoFilter()

# Update the hidden _FilterDatabase defined name to match the new range
for ($i = 1; $i -le $wb.Names.Count; $i++) {
    $n = $wb.Names.Item($i)
    if ($n.Name -like "*_FilterDatabase*") {
        $n.RefersTo = "=StoneSelection!`$A`$1:`$E`$1"
    }
}

# ---------------------------------------------------------------
# 6. Selection moves to E11 (matches the saved view state)
# ---------------------------------------------------------------
$ws.Range("E11").Select()
